$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D ("Abstand bei uns"), shifting the old
# D (Groesse im Verhaeltnis zur Erde) and E (Besonderheit) columns right.
$ws.Columns("D").Insert()

# Header for the new column.
$ws.Range("D3").Value = "Abstand bei uns"

# Fill in the new "Abstand bei uns" values for the planets that have them.
$ws.Range("D6").Value = 80
$ws.Range("D7").Value = 149
$ws.Range("D8").Value = 227
$ws.Range("D9").Value = 300
$ws.Range("D10").Value = 400
$ws.Range("D11").Value = 520

# Match the final selection left by the author.
$ws.Range("D11").Select() | Out-Null
